$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_7_7_0"
$ws.Cells.Item(2, 2).Value = 0.8918959208772095
$ws.Cells.Item(2, 3).Value = 0.9436347706764576
$ws.Cells.Item(2, 4).Value = 0.7026594591551945
$ws.Cells.Item(2, 5).Value = 0.8875413156773473
$ws.Cells.Item(2, 6).Value = 0.1196393594145775
$ws.Cells.Item(2, 7).Value = 0.07973439246416092
$ws.Cells.Item(2, 8).Value = 0.176259309053421
$ws.Cells.Item(2, 9).Value = 0.1251578778028488

$ws.Cells.Item(3, 1).Value = "model_7_7_1"
$ws.Cells.Item(3, 2).Value = 0.8945421087174573
$ws.Cells.Item(3, 3).Value = 0.9420750532697615
$ws.Cells.Item(3, 4).Value = 0.7033911405424291
$ws.Cells.Item(3, 5).Value = 0.8866751517386627
$ws.Cells.Item(3, 6).Value = 0.1167108193039894
$ws.Cells.Item(3, 7).Value = 0.08194077014923096
$ws.Cells.Item(3, 8).Value = 0.1758255660533905
$ws.Cells.Item(3, 9).Value = 0.1261218637228012

$ws.Cells.Item(4, 1).Value = "model_7_7_2"
$ws.Cells.Item(4, 2).Value = 0.8968493003884969
$ws.Cells.Item(4, 3).Value = 0.9397534814428138
$ws.Cells.Item(4, 4).Value = 0.7034390550106804
$ws.Cells.Item(4, 5).Value = 0.8851249161725033
$ws.Cells.Item(4, 6).Value = 0.114157423377037
$ws.Cells.Item(4, 7).Value = 0.08522487431764603
$ws.Cells.Item(4, 8).Value = 0.175797164440155
$ws.Cells.Item(4, 9).Value = 0.1278471499681473

$ws.Cells.Item(5, 1).Value = "model_7_7_3"
$ws.Cells.Item(5, 2).Value = 0.8993093297008945
$ws.Cells.Item(5, 3).Value = 0.938576291656831
$ws.Cells.Item(5, 4).Value = 0.7035429170334742
$ws.Cells.Item(5, 5).Value = 0.8843588224172239
$ws.Cells.Item(5, 6).Value = 0.1114348992705345
$ws.Cells.Item(5, 7).Value = 0.08689013123512268
$ws.Cells.Item(5, 8).Value = 0.1757355928421021
$ws.Cells.Item(5, 9).Value = 0.1286997497081757

$ws.Cells.Item(6, 1).Value = "model_7_7_4"
$ws.Cells.Item(6, 2).Value = 0.9016092113084652
$ws.Cells.Item(6, 3).Value = 0.9372914306333401
$ws.Cells.Item(6, 4).Value = 0.703248910312342
$ws.Cells.Item(6, 5).Value = 0.8834205135775581
$ws.Cells.Item(6, 6).Value = 0.1088896170258522
$ws.Cells.Item(6, 7).Value = 0.08870768547058105
$ws.Cells.Item(6, 8).Value = 0.1759098768234253
$ws.Cells.Item(6, 9).Value = 0.1297440081834793

$ws.Cells.Item(7, 1).Value = "model_7_7_5"
$ws.Cells.Item(7, 2).Value = 0.9038261011775033
$ws.Cells.Item(7, 3).Value = 0.936201902974387
$ws.Cells.Item(7, 4).Value = 0.7026831868179882
$ws.Cells.Item(7, 5).Value = 0.8825455533057907
$ws.Cells.Item(7, 6).Value = 0.1064361706376076
$ws.Cells.Item(7, 7).Value = 0.0902489423751831
$ws.Cells.Item(7, 8).Value = 0.1762452274560928
$ws.Cells.Item(7, 9).Value = 0.1307177841663361

$ws.Cells.Item(8, 1).Value = "model_7_7_6"
$ws.Cells.Item(8, 2).Value = 0.9058216671437656
$ws.Cells.Item(8, 3).Value = 0.9348258284816582
$ws.Cells.Item(8, 4).Value = 0.701429143210206
$ws.Cells.Item(8, 5).Value = 0.8813052346844303
$ws.Cells.Item(8, 6).Value = 0.1042276620864868
$ws.Cells.Item(8, 7).Value = 0.0921955406665802
$ws.Cells.Item(8, 8).Value = 0.1769886016845703
$ws.Cells.Item(8, 9).Value = 0.1320981532335281

$ws.Cells.Item(9, 1).Value = "model_7_7_7"
$ws.Cells.Item(9, 2).Value = 0.9077719669664925
$ws.Cells.Item(9, 3).Value = 0.9336590657903605
$ws.Cells.Item(9, 4).Value = 0.7001971352846965
$ws.Cells.Item(9, 5).Value = 0.8802112902655979
$ws.Cells.Item(9, 6).Value = 0.1020692586898804
$ws.Cells.Item(9, 7).Value = 0.09384604543447495
$ws.Cells.Item(9, 8).Value = 0.177718922495842
$ws.Cells.Item(9, 9).Value = 0.1333156228065491

$ws.Cells.Item(10, 1).Value = "model_7_7_8"
$ws.Cells.Item(10, 2).Value = 0.9096541170584966
$ws.Cells.Item(10, 3).Value = 0.9325723051436116
$ws.Cells.Item(10, 4).Value = 0.6989371819938462
$ws.Cells.Item(10, 5).Value = 0.8791641431602895
$ws.Cells.Item(10, 6).Value = 0.09998626261949539
$ws.Cells.Item(10, 7).Value = 0.09538337588310242
$ws.Cells.Item(10, 8).Value = 0.1784658133983612
$ws.Cells.Item(10, 9).Value = 0.1344810426235199

$ws.Cells.Item(11, 1).Value = "model_7_7_9"
$ws.Cells.Item(11, 2).Value = 0.9112659831305571
$ws.Cells.Item(11, 3).Value = 0.9309759331081484
$ws.Cells.Item(11, 4).Value = 0.6966770623601091
$ws.Cells.Item(11, 5).Value = 0.8775234023688564
$ws.Cells.Item(11, 6).Value = 0.0982024148106575
$ws.Cells.Item(11, 7).Value = 0.09764161705970764
$ws.Cells.Item(11, 8).Value = 0.1798055768013
$ws.Cells.Item(11, 9).Value = 0.1363070607185364

$ws.Cells.Item(12, 1).Value = "model_7_7_10"
$ws.Cells.Item(12, 2).Value = 0.9127927409562443
$ws.Cells.Item(12, 3).Value = 0.9293629743311186
$ws.Cells.Item(12, 4).Value = 0.6944423373834736
$ws.Cells.Item(12, 5).Value = 0.8758779442148074
$ws.Cells.Item(12, 6).Value = 0.09651274234056473
$ws.Cells.Item(12, 7).Value = 0.09992331266403198
$ws.Cells.Item(12, 8).Value = 0.1811303049325943
$ws.Cells.Item(12, 9).Value = 0.1381383091211319

$ws.Cells.Item(13, 1).Value = "model_7_7_11"
$ws.Cells.Item(13, 2).Value = 0.9142569736634609
$ws.Cells.Item(13, 3).Value = 0.9277867247975905
$ws.Cells.Item(13, 4).Value = 0.6924078024190639
$ws.Cells.Item(13, 5).Value = 0.87430726191187
$ws.Cells.Item(13, 6).Value = 0.09489226341247559
$ws.Cells.Item(13, 7).Value = 0.1021530777215958
$ws.Cells.Item(13, 8).Value = 0.1823363304138184
$ws.Cells.Item(13, 9).Value = 0.1398863792419434

$ws.Cells.Item(14, 1).Value = "model_7_7_12"
$ws.Cells.Item(14, 2).Value = 0.9157169111930541
$ws.Cells.Item(14, 3).Value = 0.9263873385672581
$ws.Cells.Item(14, 4).Value = 0.6907796259220935
$ws.Cells.Item(14, 5).Value = 0.8729574840064687
$ws.Cells.Item(14, 6).Value = 0.09327654540538788
$ws.Cells.Item(14, 7).Value = 0.1041326522827148
$ws.Cells.Item(14, 8).Value = 0.1833014935255051
$ws.Cells.Item(14, 9).Value = 0.1413885653018951

$ws.Cells.Item(15, 1).Value = "model_7_7_13"
$ws.Cells.Item(15, 2).Value = 0.9171600219669291
$ws.Cells.Item(15, 3).Value = 0.9251270921954675
$ws.Cells.Item(15, 4).Value = 0.6895468178923195
$ws.Cells.Item(15, 5).Value = 0.8718004499955151
$ws.Cells.Item(15, 6).Value = 0.09167943894863129
$ws.Cells.Item(15, 7).Value = 0.1059153974056244
$ws.Cells.Item(15, 8).Value = 0.1840322911739349
$ws.Cells.Item(15, 9).Value = 0.1426762491464615

$ws.Cells.Item(16, 1).Value = "model_7_7_14"
$ws.Cells.Item(16, 2).Value = 0.9182188621352736
$ws.Cells.Item(16, 3).Value = 0.9230178168299983
$ws.Cells.Item(16, 4).Value = 0.6865666308039853
$ws.Cells.Item(16, 5).Value = 0.8696339358468327
$ws.Cells.Item(16, 6).Value = 0.0905076265335083
$ws.Cells.Item(16, 7).Value = 0.1088991910219193
$ws.Cells.Item(16, 8).Value = 0.1857988834381104
$ws.Cells.Item(16, 9).Value = 0.1450874209403992

$ws.Cells.Item(17, 1).Value = "model_7_7_15"
$ws.Cells.Item(17, 2).Value = 0.9187809289822754
$ws.Cells.Item(17, 3).Value = 0.9198104402895204
$ws.Cells.Item(17, 4).Value = 0.6811468390842788
$ws.Cells.Item(17, 5).Value = 0.8661172663213549
$ws.Cells.Item(17, 6).Value = 0.08988557755947113
$ws.Cells.Item(17, 7).Value = 0.1134363412857056
$ws.Cells.Item(17, 8).Value = 0.1890116631984711
$ws.Cells.Item(17, 9).Value = 0.1490012109279633

$ws.Cells.Item(18, 1).Value = "model_7_7_16"
$ws.Cells.Item(18, 2).Value = 0.9198115212860057
$ws.Cells.Item(18, 3).Value = 0.9179098562296815
$ws.Cells.Item(18, 4).Value = 0.6788718788886903
$ws.Cells.Item(18, 5).Value = 0.8642681468094878
$ws.Cells.Item(18, 6).Value = 0.08874501287937164
$ws.Cells.Item(18, 7).Value = 0.1161249056458473
$ws.Cells.Item(18, 8).Value = 0.1903602480888367
$ws.Cells.Item(18, 9).Value = 0.1510591357946396

$ws.Cells.Item(19, 1).Value = "model_7_7_17"
$ws.Cells.Item(19, 2).Value = 0.9210418839836834
$ws.Cells.Item(19, 3).Value = 0.9166298922518741
$ws.Cells.Item(19, 4).Value = 0.6783813229428798
$ws.Cells.Item(19, 5).Value = 0.8632838021224601
$ws.Cells.Item(19, 6).Value = 0.08738337457180023
$ws.Cells.Item(19, 7).Value = 0.1179355680942535
$ws.Cells.Item(19, 8).Value = 0.1906510293483734
$ws.Cells.Item(19, 9).Value = 0.1521546244621277

$ws.Cells.Item(20, 1).Value = "model_7_7_18"
$ws.Cells.Item(20, 2).Value = 0.9220892374510219
$ws.Cells.Item(20, 3).Value = 0.9150593334570721
$ws.Cells.Item(20, 4).Value = 0.677076565599117
$ws.Cells.Item(20, 5).Value = 0.8618999191793701
$ws.Cells.Item(20, 6).Value = 0.0862242579460144
$ws.Cells.Item(20, 7).Value = 0.120157279074192
$ws.Cells.Item(20, 8).Value = 0.1914244741201401
$ws.Cells.Item(20, 9).Value = 0.1536947786808014

$ws.Cells.Item(21, 1).Value = "model_7_7_19"
$ws.Cells.Item(21, 2).Value = 0.92312064851371
$ws.Cells.Item(21, 3).Value = 0.9135850319897105
$ws.Cells.Item(21, 4).Value = 0.6761178210317957
$ws.Cells.Item(21, 5).Value = 0.8606675049594198
$ws.Cells.Item(21, 6).Value = 0.08508278429508209
$ws.Cells.Item(21, 7).Value = 0.1222428232431412
$ws.Cells.Item(21, 8).Value = 0.1919928044080734
$ws.Cells.Item(21, 9).Value = 0.1550663709640503

$ws.Cells.Item(22, 1).Value = "model_7_7_20"
$ws.Cells.Item(22, 2).Value = 0.9240800223847938
$ws.Cells.Item(22, 3).Value = 0.9120605606667429
$ws.Cells.Item(22, 4).Value = 0.6750102137189213
$ws.Cells.Item(22, 5).Value = 0.8593639855984079
$ws.Cells.Item(22, 6).Value = 0.08402103930711746
$ws.Cells.Item(22, 7).Value = 0.1243993565440178
$ws.Cells.Item(22, 8).Value = 0.1926493793725967
$ws.Cells.Item(22, 9).Value = 0.1565170884132385

$ws.Cells.Item(23, 1).Value = "model_7_7_21"
$ws.Cells.Item(23, 2).Value = 0.9249817157607654
$ws.Cells.Item(23, 3).Value = 0.9105356623950742
$ws.Cells.Item(23, 4).Value = 0.6739147903887769
$ws.Cells.Item(23, 5).Value = 0.8580632999010054
$ws.Cells.Item(23, 6).Value = 0.08302313089370728
$ws.Cells.Item(23, 7).Value = 0.126556470990181
$ws.Cells.Item(23, 8).Value = 0.1932987421751022
$ws.Cells.Item(23, 9).Value = 0.1579646319150925

$ws.Cells.Item(24, 1).Value = "model_7_7_22"
$ws.Cells.Item(24, 2).Value = 0.9258695213638402
$ws.Cells.Item(24, 3).Value = 0.9091092392944486
$ws.Cells.Item(24, 4).Value = 0.673048136538406
$ws.Cells.Item(24, 5).Value = 0.8568862581694545
$ws.Cells.Item(24, 6).Value = 0.08204060047864914
$ws.Cells.Item(24, 7).Value = 0.1285742968320847
$ws.Cells.Item(24, 8).Value = 0.1938124746084213
$ws.Cells.Item(24, 9).Value = 0.1592746078968048

$ws.Cells.Item(25, 1).Value = "model_7_7_23"
$ws.Cells.Item(25, 2).Value = 0.9266929273938589
$ws.Cells.Item(25, 3).Value = 0.9076465134111029
$ws.Cells.Item(25, 4).Value = 0.672056066264813
$ws.Cells.Item(25, 5).Value = 0.8556531732707424
$ws.Cells.Item(25, 6).Value = 0.08112932741641998
$ws.Cells.Item(25, 7).Value = 0.1306434720754623
$ws.Cells.Item(25, 8).Value = 0.1944005787372589
$ws.Cells.Item(25, 9).Value = 0.160646915435791

$ws.Cells.Item(26, 1).Value = "model_7_7_24"
$ws.Cells.Item(26, 2).Value = 0.9274771699507821
$ws.Cells.Item(26, 3).Value = 0.90621915148419
$ws.Cells.Item(26, 4).Value = 0.6710730167230894
$ws.Cells.Item(26, 5).Value = 0.8544464060926414
$ws.Cells.Item(26, 6).Value = 0.08026140183210373
$ws.Cells.Item(26, 7).Value = 0.1326626241207123
$ws.Cells.Item(26, 8).Value = 0.1949833035469055
$ws.Cells.Item(26, 9).Value = 0.1619899868965149
